$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.250616431236267
$ws.Range("B1").Value = 1.884700298309326
$ws.Range("C1").Value = 3.040507793426514
$ws.Range("D1").Value = 5.237596035003662
$ws.Range("E1").Value = 2.931501626968384
